$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price/Volume columns to Text before writing, so numeric-looking
# strings (e.g. "350.12") are stored as text rather than auto-converted numbers.
$priceVolRange = $ws.Range("D2:E51")
$priceVolRange.NumberFormat = "@"

$ws.Range("D2").Value = "51.528.95"
$ws.Range("E2").Value = "  -0.86%  "

$ws.Range("D3").Value = "2.778.55"
$ws.Range("E3").Value = "  -0.34%  "

$ws.Range("D5").Value = "350.12"
$ws.Range("E5").Value = "  -2.38%  "

$ws.Range("D6").Value = "108.15"
$ws.Range("E6").Value = "  -1.62%  "

$ws.Range("D7").Value = "0.550"
$ws.Range("E7").Value = "  -1.71%  "

$ws.Range("D8").Value = "1.00"

$ws.Range("D9").Value = "0.611"
$ws.Range("E9").Value = "  +3.53%  "

$ws.Range("D10").Value = "39.16"
$ws.Range("E10").Value = "  -1.58%  "

$ws.Range("D11").Value = "0.135"
$ws.Range("E11").Value = "  +1.57%  "

$ws.Range("D12").Value = "0.0832"
$ws.Range("E12").Value = "  -1.77%  "

$ws.Range("D13").Value = "19.82"
$ws.Range("E13").Value = "  +2.30%  "

$ws.Range("D14").Value = "7.74"
$ws.Range("E14").Value = "  +2.90%  "

$ws.Range("D15").Value = "3.221.70"
$ws.Range("E15").Value = "  -0.03%  "

$ws.Range("D16").Value = "2.784.00"
$ws.Range("E16").Value = "  -0.21%  "

$ws.Range("D17").Value = "0.923"
$ws.Range("E17").Value = "  -1.69%  "

$ws.Range("D18").Value = "51.469.22"
$ws.Range("E18").Value = "  -0.87%  "

$ws.Range("D19").Value = "7.80"
$ws.Range("E19").Value = "  +4.47%  "

$ws.Range("D20").Value = "3.08"
$ws.Range("E20").Value = "  -0.44%  "

$ws.Range("D21").Value = "13.24"
$ws.Range("E21").Value = "  +1.34%  "

$ws.Range("D22").Value = "0.0₃0962"
$ws.Range("E22").Value = "  -1.46%  "

$ws.Range("D23").Value = "69.95"
$ws.Range("E23").Value = "  -0.15%  "

$ws.Range("D24").Value = "265.92"
$ws.Range("E24").Value = "  -1.35%  "

$ws.Range("D25").Value = "2.74"
$ws.Range("E25").Value = "  -0.33%  "

$ws.Range("E26").Value = "  -0.15%  "

$ws.Range("D27").Value = "25.83"
$ws.Range("E27").Value = "  -2.44%  "

$ws.Range("D28").Value = "0.164"
$ws.Range("E28").Value = "  +0.39%  "

$ws.Range("D29").Value = "10.25"
$ws.Range("E29").Value = "  +0.35%  "

$ws.Range("D30").Value = "37.14"
$ws.Range("E30").Value = "  +9.92%  "

$ws.Range("D31").Value = "2.23"
$ws.Range("E31").Value = "  +3.97%  "

$ws.Range("D32").Value = "6.16"
$ws.Range("E32").Value = "  +7.16%  "

$ws.Range("D33").Value = "52.04"
$ws.Range("E33").Value = "  +0.22%  "

$ws.Range("D34").Value = "0.0445"
$ws.Range("E34").Value = "  -5.27%  "

$ws.Range("D35").Value = "5.55"

$ws.Range("E36").Value = "  +0.01%  "

$ws.Range("D37").Value = "0.0835"
$ws.Range("E37").Value = "  -0.48%  "

$ws.Range("D38").Value = "18.64"
$ws.Range("E38").Value = "  -0.69%  "

$ws.Range("D39").Value = "3.09"
$ws.Range("E39").Value = "  -3.34%  "

$ws.Range("D40").Value = "1.95"
$ws.Range("E40").Value = "  -2.00%  "

$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").Value = "0.114"
$ws.Range("E41").Value = "  -0.87%  "

$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "2.50"
$ws.Range("E42").Value = "  -1.82%  "

$ws.Range("D43").Value = "120.20"
$ws.Range("E43").Value = "  +0.78%  "

$ws.Range("D44").Value = "22.07"
$ws.Range("E44").Value = "  +1.41%  "

$ws.Range("E45").Value = "  -2.19%  "

$ws.Range("D46").Value = "2.140.37"
$ws.Range("E46").Value = "  +3.03%  "

$ws.Range("D47").Value = "3.28"
$ws.Range("E47").Value = "  +1.36%  "

$ws.Range("D48").Value = "2.33"
$ws.Range("E48").Value = "  +5.03%  "

$ws.Range("D49").Value = "0.227"
$ws.Range("E49").Value = "  +19.51%  "

$ws.Range("D50").Value = "5.44"
$ws.Range("E50").Value = "  -5.06%  "

$ws.Range("D51").Value = "0.903"
$ws.Range("E51").Value = "  -5.29%  "

# Restore original (default) formatting so styling matches source
$priceVolRange.Style = "Normal"